$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells in row 1: drop "Checkbox" from "input_rowSelectionCheckbox_*"
for ($col = 1; $col -le 12; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = $cell.Value2 -replace "input_rowSelectionCheckbox_", "input_rowSelection_"
}

# Narrow column L from 50 to 42 (raw OOXML width units).
# The ColumnWidth COM property is offset from the stored OOXML width by the
# sheet's default character-padding (~0.83 here), so compensate accordingly.
$ws.Columns.Item(12).ColumnWidth = 41.17
